$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "Test12bcv usernamer"
$ws.Range("B3").Value = "PrbfjaEqwkil1111@gmail.com"
$ws.Range("A4").Value = "Testdf usern!_2321"
$ws.Range("B4").Value = "Proper11127451@gmail.com"
$ws.Range("A5").Value = "Test usernameßüöäe"
$ws.Range("B5").Value = "Pr2operEmaila1222!2_2@gmail.com"

$ws.Range("A5").Select()
